$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.339.50'
$ws.Range('E2').Value = '  -2.31%  '
$ws.Range('D3').Value = '2.580.42'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.37'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.11'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('D9').Value = '2.587.96'
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('E12').Value = '  +12.54%  '
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('D14').Value = '3.035.35'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').Value = '59.302.72'
$ws.Range('E15').Value = '  -2.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.98'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +4.14%  '
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '2.581.39'
$ws.Range('E18').Value = '  -2.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.55'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.65'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.35'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.47'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.476'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +8.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '62.54'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -5.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.42'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = '0.0₃0776'
$ws.Range('E29').Value = '  -3.31%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.25'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.86'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.09'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.11'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  +1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.899'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.40'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.854'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -4.36%  '
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.68'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '289.84'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -3.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.31'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +6.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0974'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.594'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0531'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0235'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.73'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').Value = '1.942.20'
$ws.Range('E51').Value = '  -0.95%  '
